$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings (e.g. "228.01",
# "1.514.07") are preserved as text instead of being parsed into numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '40.041.75'
$ws.Range('E2').Value = '  +1.44%  '

$ws.Range('D3').Value = '2.194.03'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '228.01'
$ws.Range('E5').Value = '  -0.45%  '

$ws.Range('D6').Value = '0.628'
$ws.Range('E6').Value = '  +1.01%  '

$ws.Range('D7').Value = '63.25'
$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  -0.61%  '

$ws.Range('D10').Value = '0.0861'
$ws.Range('E10').Value = '  -0.81%  '

$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  -0.08%  '

$ws.Range('D12').Value = '2.517.78'
$ws.Range('E12').Value = '  +1.32%  '

$ws.Range('D13').Value = '15.80'
$ws.Range('E13').Value = '  -1.42%  '

$ws.Range('D14').Value = '21.98'
$ws.Range('E14').Value = '  -1.06%  '

$ws.Range('D15').Value = '0.815'
$ws.Range('E15').Value = '  -0.29%  '

$ws.Range('E16').Value = '  -0.60%  '

$ws.Range('D17').Value = '2.191.92'
$ws.Range('E17').Value = '  +1.51%  '

$ws.Range('D18').Value = '39.921.50'
$ws.Range('E18').Value = '  +1.11%  '

$ws.Range('D19').Value = '0.0₃0906'
$ws.Range('E19').Value = '  +6.03%  '

$ws.Range('D20').Value = '72.16'
$ws.Range('E20').Value = '  -0.12%  '

$ws.Range('E21').Value = '  -1.61%  '

$ws.Range('D22').Value = '232.40'
$ws.Range('E22').Value = '  +1.37%  '

$ws.Range('E23').Value = '  +0.08%  '

$ws.Range('D24').Value = '2.34'
$ws.Range('E24').Value = '  -0.97%  '

$ws.Range('E25').Value = '  +0.42%  '

$ws.Range('D26').Value = '9.65'
$ws.Range('E26').Value = '  -1.12%  '

$ws.Range('D27').Value = '171.50'
$ws.Range('E27').Value = '  -0.49%  '

$ws.Range('D28').Value = '0.141'
$ws.Range('E28').Value = '  +2.25%  '

$ws.Range('E29').Value = '  +2.31%  '

$ws.Range('D30').Value = '20.09'
$ws.Range('E30').Value = '  +1.85%  '

$ws.Range('D31').Value = '2.74'
$ws.Range('E31').Value = '  +4.38%  '

$ws.Range('E32').Value = '  +0.49%  '

$ws.Range('D33').Value = '4.57'
$ws.Range('E33').Value = '  -2.07%  '

$ws.Range('E34').Value = '  -2.56%  '

$ws.Range('D35').Value = '6.97'
$ws.Range('E35').Value = '  -1.57%  '

$ws.Range('D36').Value = '0.0624'
$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('D37').Value = '3.87'
$ws.Range('E37').Value = '  +6.51%  '

$ws.Range('E38').Value = '  -0.08%  '

$ws.Range('D39').Value = '5.17'
$ws.Range('E39').Value = '  +22.17%  '

$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.03%  '

$ws.Range('D41').Value = '102.96'
$ws.Range('E41').Value = '  -1.77%  '

$ws.Range('E42').Value = '  -0.93%  '

$ws.Range('D43').Value = '1.23'
$ws.Range('E43').Value = '  +2.00%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '17.41'
$ws.Range('E44').Value = '  -2.02%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.514.07'
$ws.Range('E45').Value = '  -1.62%  '

$ws.Range('D46').Value = '8.26'
$ws.Range('E46').Value = '  +4.56%  '

$ws.Range('E47').Value = '  -0.78%  '

$ws.Range('E48').Value = '  -0.78%  '

$ws.Range('E49').Value = '  -0.28%  '

$ws.Range('D50').Value = '0.000197'
$ws.Range('E50').Value = '  +33.59%  '

$ws.Range('D51').Value = '50.17'
$ws.Range('E51').Value = '  +7.29%  '

# Restore the default (Normal) style so the cell styling matches the original workbook.
$dataRange.Style = "Normal"
